# Update host species labels (use common name + scientific name) on the
# "Antibodies" sheet, add more host species / the IgY isotype on the
# "Terminology" reference sheet, and widen the data-validation ranges that
# point at the newly-grown Terminology lists.

$wb = $excel.ActiveWorkbook

$antibodies = $wb.Worksheets.Item("Antibodies")
$terminology = $wb.Worksheets.Item("Terminology")

# ---------------------------------------------------------------------------
# "Antibodies" sheet - column B ("Host") - expand the short Latin names to
# "<common name> (<Latin name>)", matching the new Terminology vocabulary.
# ---------------------------------------------------------------------------
$antibodies.Range("B2").Value = "human (Homo sapiens)"
$antibodies.Range("B3").Value = "human (Homo sapiens)"
$antibodies.Range("B4").Value = "mouse (Mus musculus)"
$antibodies.Range("B6").Value = "mouse (Mus musculus)"
$antibodies.Range("B7").Value = "mouse (Mus musculus)"
$antibodies.Range("B8").Value = "human (Homo sapiens)"
$antibodies.Range("B9").Value = "mouse (Mus musculus)"
$antibodies.Range("B10").Value = "human (Homo sapiens)"
$antibodies.Range("B11").Value = "mouse (Mus musculus)"

# Data validation lists need to cover the grown Terminology ranges: Host
# grew from 3 rows (A2:A4) to 5 rows (A2:A6), Isotype grew by one row
# (B2:B15 -> B2:B16) to fit the new "IgY" entry.
$antibodies.Range("B2:B100").Validation.Formula1 = "=Terminology!`$A`$2:`$A`$6"
$antibodies.Range("C2:C100").Validation.Formula1 = "=Terminology!`$B`$2:`$B`$16"

# ---------------------------------------------------------------------------
# "Terminology" sheet - protected, so unprotect, edit, then re-protect.
# Column A ("Host") now lists five species with common names, and column B
# ("Isotype") gains "IgY" in row 16.
# ---------------------------------------------------------------------------
$wasProtected = $terminology.ProtectContents
if ($wasProtected) {
    $terminology.Unprotect()
}

$terminology.Range("A2").Value = "chicken (Gallus gallus)"
$terminology.Range("A3").Value = "human (Homo sapiens)"
$terminology.Range("A4").Value = "llama (Lama glama)"
$terminology.Range("A5").Value = "mouse (Mus musculus)"
$terminology.Range("A6").Value = "alpaca (Vicugna pacos)"
$terminology.Range("B16").Value = "IgY"

if ($wasProtected) {
    $terminology.Protect()
}
